$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: remove the "RUN" marker from A7 (moved down to the new row 8) ---
$ws.Range("A7").ClearContents()

# --- Row 8: populate with a new test-run entry (copy of row 7's pattern) ---
$ws.Range("A8").Value = "run"
$ws.Range("B8").Value = "DPLKAKT130-001"
$ws.Range("C8").Value = "Akuntansi - Transaksi"
$ws.Range("D8").Value = "Cek Jurnal Settlement Transaksi Klaim P01"
$ws.Range("E8").Value = "Jurnal Settlement Terbentuk dengan Benar di-DPLK"
$ws.Range("G8").Value = 34786
$ws.Range("H8").Value = "bni1234"
$ws.Range("I8").Value = "Keuangan"
$ws.Range("J8").Value = "Transaksi"
$ws.Range("K8").Value = "Keuangan Kepesertaan"
$ws.Range("L8").Value = "Inquiry Pembayaran Kepesertaan"
$ws.Range("M8").Value = "DTOBL202300007"
$ws.Range("P8").Value = "BC001"

$f8 = "=""Username : ""&G8&"",`nPassword : ""&H8&"",`nDokumen ID : ""&M8&"",`nTanggal : ""&N8"
$ws.Range("F8").Formula = $f8

# J8 carries a quote-prefixed style (xf 3) in the template row; restore it
# since writing .Value above reset it back to the column's default style.
$ws.Range("J7").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: matches where the user left the selection after editing ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D6").Select()
